# Notas do forum para a semana 05/06/2022 a 11/06/2022 no semestre 20022-1
#
# The sheet has per-day view columns B..J (2022-05-29 .. 2022-06-06) followed
# by summary columns K (total_views) and L (nota_view). This week's refresh
# drops the two oldest day columns (I = 2022-06-05, J = 2022-06-06) and
# recomputes the summary so it lands back in I:J, shrinking the sheet from
# A1:L47 to A1:J47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the two old date columns (I:J) shifts the summary columns
# (old K -> new I "total_views", old L -> new J "nota_view") left for every
# row, including the header text in row 1.
$ws.Range("I1:J1").EntireColumn.Delete()

# A handful of students had nonzero activity on the two dropped days, so
# their recomputed total_views/nota_view differ from a pure left-shift of
# the old summary columns. Patch those rows to the recomputed values.
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 2

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 2

$ws.Range("I12").Value = 4
$ws.Range("J12").Value = 2

$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 2

$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0

$ws.Range("I23").Value = 5
$ws.Range("J23").Value = 2

$ws.Range("I35").Value = 7
$ws.Range("J35").Value = 2

$ws.Range("I41").Value = 1
$ws.Range("J41").Value = 0.5

$ws.Range("I43").Value = 6
$ws.Range("J43").Value = 2

$ws.Range("I44").Value = 4
$ws.Range("J44").Value = 2
